# Added frontend and .gitignore file
# - Add a new "Frontend" worksheet at the end of the workbook with the
#   Maven frontend setup info, and make it the active sheet/tab.
# - Leave the selection on Sheet1 parked at A11 (already saved in the
#   source file) but it is no longer the active sheet, so Excel drops its
#   tabSelected flag once "Frontend" becomes active.

$wb = $excel.ActiveWorkbook

# Scroll Sheet1's window so row 11 is pinned at the top (topLeftCell="A11")
# while keeping the existing A11 selection untouched.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

# Grab the current last sheet ("Change Sets") so the new sheet is appended
# after it, matching the sheet order in the diff (Sheet1, Change Sets, Frontend).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$frontend = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$frontend.Name = "Frontend"

# Populate the new sheet's header row.
$frontend.Range("A1").Value = "Setup Info"
$frontend.Range("B1").Value = "Maven Frontend Goal"

# Land the selection on B2 (as if the user had just tabbed/entered through
# the two header cells) and make sure Frontend is the active/selected sheet.
$frontend.Range("B2").Select() | Out-Null
$frontend.Activate() | Out-Null
